$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.246.11"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.862.25"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7011"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.32"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07658"
$ws.Range("E8").Value = "  +0.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3053"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.30"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08173"
$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("D12").Value = "1.865.84"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7179"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("E14").Value = "  -0.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.62"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").Value = "29.239.19"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.748"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.04"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +0.60%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007705"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "2.110.91"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.486"
$ws.Range("E24").Value = "  -1.66%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.30"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.016"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1458"
$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.06"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.973"
$ws.Range("E29").Value = "  +2.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.410"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.453"
$ws.Range("E31").Value = "  +0.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.480"
$ws.Range("E32").Value = "  -1.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.000"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05190"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.164"
$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7078"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.656"
$ws.Range("E38").Value = "  -0.33%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.718"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9333"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").Value = "1.141.59"
$ws.Range("E42").Value = "  +9.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4283"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.79"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.873"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.35"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.805"
$ws.Range("E48").Value = "  +3.92%  "

$ws.Range("D49").Value = "2.008.21"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.141"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.969"
$ws.Range("E51").Value = "  -3.64%  "
